$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# -----------------------------------------------------------------
# 1. Refresh the "time_taken" timestamps in column F of the "data"
#    sheet (rows 2..103) to reflect the new panel-query run.
# -----------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:33:25.424678",
    "2021-10-05 14:33:25.424683",
    "2021-10-05 14:33:25.424686",
    "2021-10-05 14:33:25.424688",
    "2021-10-05 14:33:25.424690",
    "2021-10-05 14:33:25.424692",
    "2021-10-05 14:33:25.424694",
    "2021-10-05 14:33:25.424695",
    "2021-10-05 14:33:25.424697",
    "2021-10-05 14:33:25.424699",
    "2021-10-05 14:33:25.424701",
    "2021-10-05 14:33:25.424703",
    "2021-10-05 14:33:25.424705",
    "2021-10-05 14:33:25.424707",
    "2021-10-05 14:33:25.424709",
    "2021-10-05 14:33:25.424711",
    "2021-10-05 14:33:25.424713",
    "2021-10-05 14:33:25.424714",
    "2021-10-05 14:33:25.424716",
    "2021-10-05 14:33:25.424718",
    "2021-10-05 14:33:25.424720",
    "2021-10-05 14:33:25.424722",
    "2021-10-05 14:33:25.424724",
    "2021-10-05 14:33:25.424726",
    "2021-10-05 14:33:25.424728",
    "2021-10-05 14:33:25.424730",
    "2021-10-05 14:33:25.424732",
    "2021-10-05 14:33:25.424734",
    "2021-10-05 14:33:25.424735",
    "2021-10-05 14:33:25.424737",
    "2021-10-05 14:33:25.424739",
    "2021-10-05 14:33:25.424741",
    "2021-10-05 14:33:25.424743",
    "2021-10-05 14:33:25.424745",
    "2021-10-05 14:33:25.424747",
    "2021-10-05 14:33:25.424749",
    "2021-10-05 14:33:25.424751",
    "2021-10-05 14:33:25.424753",
    "2021-10-05 14:33:25.424755",
    "2021-10-05 14:33:25.424757",
    "2021-10-05 14:33:25.424759",
    "2021-10-05 14:33:25.424761",
    "2021-10-05 14:33:25.424763",
    "2021-10-05 14:33:25.424764",
    "2021-10-05 14:33:25.424767",
    "2021-10-05 14:33:25.424768",
    "2021-10-05 14:33:25.424770",
    "2021-10-05 14:33:25.424772",
    "2021-10-05 14:33:25.424774",
    "2021-10-05 14:33:25.424776",
    "2021-10-05 14:33:25.424778",
    "2021-10-05 14:33:25.424780",
    "2021-10-05 14:33:25.424782",
    "2021-10-05 14:33:25.424784",
    "2021-10-05 14:33:25.424786",
    "2021-10-05 14:33:25.424788",
    "2021-10-05 14:33:25.424790",
    "2021-10-05 14:33:25.424792",
    "2021-10-05 14:33:25.424794",
    "2021-10-05 14:33:25.424796",
    "2021-10-05 14:33:25.424798",
    "2021-10-05 14:33:25.424799",
    "2021-10-05 14:33:25.424801",
    "2021-10-05 14:33:25.424803",
    "2021-10-05 14:33:25.424806",
    "2021-10-05 14:33:25.424808",
    "2021-10-05 14:33:25.424810",
    "2021-10-05 14:33:25.424812",
    "2021-10-05 14:33:25.424814",
    "2021-10-05 14:33:25.424815",
    "2021-10-05 14:33:25.424817",
    "2021-10-05 14:33:25.424819",
    "2021-10-05 14:33:25.424821",
    "2021-10-05 14:33:25.424823",
    "2021-10-05 14:33:25.424825",
    "2021-10-05 14:33:25.424827",
    "2021-10-05 14:33:25.424830",
    "2021-10-05 14:33:25.424833",
    "2021-10-05 14:33:25.424835",
    "2021-10-05 14:33:25.424837",
    "2021-10-05 14:33:25.424839",
    "2021-10-05 14:33:25.424841",
    "2021-10-05 14:33:25.424843",
    "2021-10-05 14:33:25.424845",
    "2021-10-05 14:33:25.424847",
    "2021-10-05 14:33:25.424849",
    "2021-10-05 14:33:25.424851",
    "2021-10-05 14:33:25.424852",
    "2021-10-05 14:33:25.424855",
    "2021-10-05 14:33:25.424857",
    "2021-10-05 14:33:25.424859",
    "2021-10-05 14:33:25.424861",
    "2021-10-05 14:33:25.424863",
    "2021-10-05 14:33:25.424866",
    "2021-10-05 14:33:25.424868",
    "2021-10-05 14:33:25.424870",
    "2021-10-05 14:33:25.424871",
    "2021-10-05 14:33:25.424873",
    "2021-10-05 14:33:25.424875",
    "2021-10-05 14:33:25.424877",
    "2021-10-05 14:33:25.424879",
    "2021-10-05 14:33:25.424881"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $newTimes[$i]
}

# -----------------------------------------------------------------
# 2. Add a new "metadata" worksheet after "data" describing the
#    panel query that produced this workbook. Duplicate "data" first
#    so the new sheet inherits the same sheetPr/page-setup, then wipe
#    its contents — cleaner than Worksheets.Add(), which stamps a
#    generic blank sheet with different default margins.
# -----------------------------------------------------------------
$ws.Copy($null, $ws)
$meta = $wb.Worksheets.Item($ws.Index + 1)
$meta.Name = "metadata"
$meta.Cells.Clear()

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Cerebral vascular malformations"
$meta.Range("C2").Value = 3144
$meta.Range("E2").Value = "2021-09-20T02:40:05.024937Z"
$meta.Range("F2").Value = "2021-10-05 14:33:25.422213"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3144/?format=json"

# data_version ("0.23") must stay a *text* value, not be coerced into
# the number 0.23 — stage it on a scratch cell formatted as Text, then
# copy only the value across so D2 keeps the sheet's default (unstyled)
# formatting, matching the rest of the row.
$scratch = $meta.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "0.23"
$scratch.Copy()
$meta.Range("D2").PasteSpecial(-4163)
$scratch.Clear()

# Match the header row / index-column styling used on the "data" sheet
# (bold, centred, thin-bordered) by copying its format, instead of
# rebuilding it property-by-property (which would mint a near-duplicate
# style entry).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
$ws.Activate()

Write-Output "done"
